# feat: add multi rules tables
#
# Adds a new "Multi" worksheet (a rule table similar to the existing
# "Rules" sheet's Rule1 table, but for a multi-result "Multi1" rule with
# three result classes A/B/C) and updates the active-sheet/selection
# state left behind by the edit.

$wb = $excel.ActiveWorkbook

$rules = $wb.Worksheets.Item("Rules")

# --- 1. Create the new "Multi" sheet at the end of the tab strip -----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$multi = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$multi.Name = "Multi"

# --- 2. Write all cell values first (plain cells, no formatting yet) -------
# Re-used labels/values (already existing shared strings)
$multi.Range("D7").Value = "C1"
$multi.Range("E7").Value = "HC1"
$multi.Range("F7").Value = "RET1"
$multi.Range("D8").Value = "measure1"
$multi.Range("E8").Value = "measure2"
$multi.Range("D9").Value = "String"
$multi.Range("E9").Value = "IntRange"
$multi.Range("F9").Value = "String"

# New "C" result row (row 13)
$multi.Range("E13").Value = "C"
$multi.Range("F13").Value = "C"
$multi.Range("G13").Value = "C"
$multi.Range("H13").Value = "C"

# "A" result row (row 11)
$multi.Range("D11").Value = "A"
$multi.Range("E11").Value = "A1"
$multi.Range("F11").Value = "A2"
$multi.Range("G11").Value = "A3"
$multi.Range("H11").Value = "A4"

# "B" result row (row 12)
$multi.Range("D12").Value = "B"
$multi.Range("E12").Value = "B1"
$multi.Range("F12").Value = "B2"
$multi.Range("G12").Value = "B3"
$multi.Range("H12").Value = "B4"

# measure1 threshold row (row 10)
$multi.Range("D10").Value = "measure1"
$multi.Range("F10").Value = "[25..66]"
$multi.Range("G10").Value = "[33..74]"
$multi.Range("E10").Value = "<50"
$multi.Range("H10").Value = ">=75"

# Header / rule signature (row 6) - rich text: bold rule name in the middle
$header = $multi.Range("D6")
$header.Value = "Multi String Multi1(String measure1, Integer measure2)"
$boldPart = $header.Characters(14, 6)
$boldPart.Font.Bold = $true
$boldPart.Font.Color = 12566463
$restPart = $header.Characters(20, 36)
$restPart.Font.Color = 12566463

# --- 3. Apply the table's formatting (borders/fills/fonts) by pasting the --
#        formats from the equivalent "Rule1" table on the "Rules" sheet, one
#        column to the left and two rows up. Formats only - values already
#        set above are left untouched.
$rules.Range("B3:H11").Copy() | Out-Null
$multi.Range("C5").PasteSpecial(-4122) | Out-Null

$rules.Range("B11:H11").Copy() | Out-Null
$multi.Range("C14").PasteSpecial(-4122) | Out-Null

$rules.Range("B10:H10").Copy() | Out-Null
$multi.Range("C13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 4. Misc view-state cleanup left by the edit ---------------------------
# The "Rules" sheet selection moved to the header of its own table.
$rules.Activate() | Out-Null
$rules.Range("B3:H11").Select() | Out-Null

# The new sheet is the one left active/selected.
$multi.Activate() | Out-Null
$multi.Range("F15").Select() | Out-Null
